$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 11155
$ws.Range("F4").Value = 278
$ws.Range("F5").Value = 1246
$ws.Range("F6").Value = 1119
$ws.Range("F7").Value = 867
$ws.Range("F8").Value = 295
$ws.Range("F10").Value = 1194
$ws.Range("F12").Value = 927
$ws.Range("F13").Value = 2141
$ws.Range("F14").Value = 25
$ws.Range("F15").Value = 1057
$ws.Range("F16").Value = 850
$ws.Range("F17").Value = 566
$ws.Range("F18").Value = 833
$ws.Range("F19").Value = 958
$ws.Range("F21").Value = 273
$ws.Range("F23").Value = 661
$ws.Range("F24").Value = 685
$ws.Range("F25").Value = 139
$ws.Range("F26").Value = 379
$ws.Range("F28").Value = 52
$ws.Range("F29").Value = 141
$ws.Range("F30").Value = 517
$ws.Range("F31").Value = 185
$ws.Range("F32").Value = 261
$ws.Range("F33").Value = 254
$ws.Range("F34").Value = 601
$ws.Range("F35").Value = 2135
$ws.Range("F36").Value = 411
$ws.Range("F37").Value = 56
$ws.Range("F38").Value = 1475
$ws.Range("F39").Value = 415
$ws.Range("F40").Value = 127
$ws.Range("F41").Value = 59
$ws.Range("F43").Value = 49
$ws.Range("F45").Value = 87
$ws.Range("F47").Value = 54

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 204
$ws.Range("F6").Value = 14
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 4403
$ws.Range("F17").Value = 9
$ws.Range("F18").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2203
$ws.Range("F3").Value = 655
$ws.Range("F4").Value = 607

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2203
$ws.Range("F3").Value = 655
$ws.Range("F5").Value = 278
$ws.Range("F6").Value = 1246
$ws.Range("F8").Value = 607
$ws.Range("F9").Value = 1119
$ws.Range("F10").Value = 204
$ws.Range("F11").Value = 295
$ws.Range("F12").Value = 1194
$ws.Range("F13").Value = 14
$ws.Range("F16").Value = 927
$ws.Range("F17").Value = 2141
$ws.Range("F18").Value = 25
$ws.Range("F19").Value = 1057
$ws.Range("F20").Value = 850
$ws.Range("F21").Value = 566
$ws.Range("F22").Value = 833
$ws.Range("F23").Value = 958
$ws.Range("F25").Value = 661
$ws.Range("F27").Value = 685
$ws.Range("F28").Value = 139
$ws.Range("F29").Value = 379
$ws.Range("F31").Value = 52
$ws.Range("F32").Value = 517
$ws.Range("F33").Value = 185
$ws.Range("F34").Value = 261
$ws.Range("F35").Value = 254
$ws.Range("F36").Value = 2135
$ws.Range("F37").Value = 148
$ws.Range("F38").Value = 411
$ws.Range("F39").Value = 56
$ws.Range("F40").Value = 1475
$ws.Range("F41").Value = 415
$ws.Range("F42").Value = 127
$ws.Range("F43").Value = 59
$ws.Range("F46").Value = 87
$ws.Range("F48").Value = 54
